$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.031" or
# "0.000008992" are not auto-coerced to numbers by Excel's input parser.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.526.90"
$ws.Range("E2").Value = "  +4.00%  "

$ws.Range("D3").Value = "1.841.74"
$ws.Range("E3").Value = "  +2.91%  "

$ws.Range("D4").Value = "1.031"
$ws.Range("E4").Value = "  +2.94%  "

$ws.Range("D5").Value = "318.77"
$ws.Range("E5").Value = "  +4.40%  "

$ws.Range("D6").Value = "1.024"
$ws.Range("E6").Value = "  +2.23%  "

$ws.Range("D7").Value = "0.4369"
$ws.Range("E7").Value = "  +2.67%  "

$ws.Range("D8").Value = "0.3734"
$ws.Range("E8").Value = "  +3.01%  "

$ws.Range("D9").Value = "0.07375"
$ws.Range("E9").Value = "  +2.92%  "

$ws.Range("D10").Value = "0.8739"
$ws.Range("E10").Value = "  +2.45%  "

$ws.Range("E11").Value = "  +4.16%  "

$ws.Range("D12").Value = "1.847.84"
$ws.Range("E12").Value = "  +2.79%  "

$ws.Range("D13").Value = "5.493"
$ws.Range("E13").Value = "  +4.33%  "

$ws.Range("D14").Value = "6.677"
$ws.Range("E14").Value = "  +2.84%  "

$ws.Range("D15").Value = "0.07152"
$ws.Range("E15").Value = "  +3.55%  "

$ws.Range("D16").Value = "82.64"
$ws.Range("E16").Value = "  +3.86%  "

$ws.Range("D17").Value = "1.032"
$ws.Range("E17").Value = "  +2.37%  "

$ws.Range("D18").Value = "0.000008992"
$ws.Range("E18").Value = "  +2.45%  "

$ws.Range("E19").Value = "  +2.45%  "

$ws.Range("D20").Value = "15.38"
$ws.Range("E20").Value = "  +2.57%  "

$ws.Range("D21").Value = "27.520.92"
$ws.Range("E21").Value = "  +3.89%  "

$ws.Range("E22").Value = "  +2.26%  "

$ws.Range("D23").Value = "11.19"
$ws.Range("E23").Value = "  +1.35%  "

$ws.Range("D24").Value = "2.065.23"
$ws.Range("E24").Value = "  +2.24%  "

$ws.Range("D25").Value = "157.54"
$ws.Range("E25").Value = "  +3.41%  "

$ws.Range("D26").Value = "1.927"
$ws.Range("E26").Value = "  +6.02%  "

$ws.Range("E27").Value = "  +3.07%  "

$ws.Range("D28").Value = "5.249"
$ws.Range("E28").Value = "  +2.27%  "

$ws.Range("D29").Value = "1.931"
$ws.Range("E29").Value = "  +1.40%  "

$ws.Range("D30").Value = "115.84"
$ws.Range("E30").Value = "  +1.02%  "

$ws.Range("D31").Value = "0.09093"
$ws.Range("E31").Value = "  +2.22%  "

$ws.Range("D32").Value = "1.203"
$ws.Range("E32").Value = "  +5.35%  "

$ws.Range("D33").Value = "0.7666"
$ws.Range("E33").Value = "  +2.99%  "

$ws.Range("D34").Value = "4.499"
$ws.Range("E34").Value = "  +3.32%  "

$ws.Range("D35").Value = "2.877"
$ws.Range("E35").Value = "  +4.30%  "

$ws.Range("D36").Value = "1.029"
$ws.Range("E36").Value = "  +2.76%  "

$ws.Range("D37").Value = "1.144"
$ws.Range("E37").Value = "  +3.18%  "

$ws.Range("D38").Value = "0.01970"
$ws.Range("E38").Value = "  +3.77%  "

$ws.Range("D39").Value = "0.05246"
$ws.Range("E39").Value = "  +1.62%  "

$ws.Range("D40").Value = "0.5164"
$ws.Range("E40").Value = "  +3.84%  "

$ws.Range("D41").Value = "2.786"
$ws.Range("E41").Value = "  +6.69%  "

$ws.Range("D42").Value = "0.1670"
$ws.Range("E42").Value = "  +2.98%  "

$ws.Range("D43").Value = "6.654"
$ws.Range("E43").Value = "  +4.55%  "

$ws.Range("D44").Value = "8.523"
$ws.Range("E44").Value = "  +3.89%  "

$ws.Range("D45").Value = "108.77"
$ws.Range("E45").Value = "  +3.22%  "

$ws.Range("D46").Value = "10.59"
$ws.Range("E46").Value = "  +2.50%  "

$ws.Range("D47").Value = "1.713"
$ws.Range("E47").Value = "  +4.22%  "

$ws.Range("D48").Value = "0.4641"
$ws.Range("E48").Value = "  +2.76%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.06358"
$ws.Range("E49").Value = "  +2.52%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "1.889"
$ws.Range("E50").Value = "  +7.50%  "

$ws.Range("D51").Value = "39.42"
$ws.Range("E51").Value = "  +6.77%  "

# Restore the default (unstyled) cell format for column D so the saved
# workbook matches the original's styling (no explicit style index).
$ws.Range("D2:D51").Style = "Normal"
